$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actuator Calc Assumptions")

# Rewrite the loading/dimension data table (rows 2-5), adding two new rows
# for "Tread Radius" and "Distance to Center of Gravity", renaming "Mass"
# to "Total Mass", and moving "Operating Speed" to the bottom.
$ws.Range("A2").Value = "Total Mass"
$ws.Range("B2").Value = "30 000 kg"

$ws.Range("A3").Value = "Tread Radius"
$ws.Range("B3").Value = "0.5 m"

$ws.Range("A4").Value = "Distance to Center of Gravity"
$ws.Range("B4").Value = "1.92 m"

$ws.Range("A5").Value = "Operating Speed"
$ws.Range("B5").Value = "3 m/s"

# Column A now needs to fit the longer labels (e.g. "Distance to Center of
# Gravity"), so widen it to match the new best-fit width.
$ws.Columns.Item(1).ColumnWidth = 30.2

$ws.Range("A6").Select()
